# "Update plots and report"
# Refresh the Sheet1 raw data (columns K:N, rows 26-30) with newly-rescaled
# benchmark figures; the G (sum) and H (ratio) formulas in those rows pick
# up the new totals automatically through recalculation. Also restores the
# view state (scroll position / selected cell / window width) left behind
# by the author after making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Updated benchmark breakdown values -----------------------------------
# Row 26 (key size 192)
$ws.Range("K26").Value = 92322
$ws.Range("L26").Value = 199576
$ws.Range("M26").Value = 52683
$ws.Range("N26").Value = 0

# Row 27 (key size 224)
$ws.Range("K27").Value = 139226
$ws.Range("L27").Value = 256289
$ws.Range("M27").Value = 81518
$ws.Range("N27").Value = 0

# Row 28 (key size 256)
$ws.Range("K28").Value = 158049
$ws.Range("L28").Value = 310193
$ws.Range("M28").Value = 89900
$ws.Range("N28").Value = 0

# Row 29 (key size 384)
$ws.Range("K29").Value = 350800
$ws.Range("L29").Value = 619485
$ws.Range("M29").Value = 193140
$ws.Range("N29").Value = 0

# Row 30 (key size 521)
$ws.Range("K30").Value = 725781
$ws.Range("L30").Value = 1067380
$ws.Range("M30").Value = 395396
$ws.Range("N30").Value = 0

$excel.CalculateFull()

# --- View / window state ---------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("O31").Select() | Out-Null
$excel.ActiveWindow.Width = 591.25
